$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 17 for the new "능력치보기" (View Stats) place
$ws.Range("A17").Value = 10001
$ws.Range("B17").Value = "능력치보기"
$ws.Range("D17").Value = "{1003}"
$ws.Range("C17").Value = " "

# Update Selections for PlaceID 1003 ("내 집") row (row 4) to include 10001
$ws.Range("D4").Value = "{10001,1001,1002}"

# Update selection to match the author's final cursor position
$ws.Range("E6").Select()
